$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Extend the "Tabella1" table by one row (A1:F24 -> A1:F25) ---
$lo = $ws.ListObjects.Item("Tabella1")
$lo.Resize($ws.Range("A1:F25"))

# New row 25 becomes the new "last row" of the table and should carry the
# bottom-border formatting that row 24 used to have (snapshot it first).
$ws.Range("A24:F24").Copy()
$ws.Range("A25:F25").PasteSpecial(-4122)   # xlPasteFormats

# Row 24 currently still carries the "last row" (bottom-border) formatting;
# it needs to become a normal interior row, matching rows 2-23.
$ws.Range("A2:F2").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the new food item: gnocchi smart
$ws.Range("A25").Value = "gnocchi smart"
$ws.Range("B25").Value = 4.5
$ws.Range("C25").Value = 174
$ws.Range("D25").Formula = "=B25/C25"
$ws.Range("E25").Value = 1
$ws.Range("F25").Formula = "=B25/E25"

# --- Below-table filler rows: extend one further row (39 -> 40) ---
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)

$ws.Range("E39").Copy()
$ws.Range("E40").PasteSpecial(-4122)

# --- Restore the selected cell shown in the workbook ---
$ws.Range("Q10").Select()
